# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with refreshed data, matching a GitHub Actions scheduled data-refresh.
# Row 15/16 (WrappedBTC/ShibaInu) swap position with new values.
# A leading apostrophe ('' in a single-quoted literal) forces Excel to store
# numeric-looking price strings (e.g. "0.999") as text, matching the
# original inline-string cell type instead of converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.424.82'
$ws.Range("E2").Value = '  +2.05%  '

$ws.Range("D3").Value = '2.981.69'
$ws.Range("E3").Value = '  +2.35%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '''503.42'
$ws.Range("E5").Value = '  +7.52%  '

$ws.Range("D6").Value = '''134.71'
$ws.Range("E6").Value = '  +8.84%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").Value = '''0.428'
$ws.Range("E8").Value = '  +6.72%  '

$ws.Range("E9").Value = '  +11.50%  '

$ws.Range("E10").Value = '  +12.06%  '

$ws.Range("D11").Value = '''0.351'
$ws.Range("E11").Value = '  +7.47%  '

$ws.Range("E12").Value = '  +3.19%  '

$ws.Range("D13").Value = '3.491.13'
$ws.Range("E13").Value = '  +2.14%  '

$ws.Range("D14").Value = '''25.08'
$ws.Range("E14").Value = '  +12.16%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.0000151'
$ws.Range("E15").Value = '  +14.94%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '56.357.83'
$ws.Range("E16").Value = '  +1.89%  '

$ws.Range("D17").Value = '2.977.61'

$ws.Range("D18").Value = '''5.66'
$ws.Range("E18").Value = '  +11.85%  '

$ws.Range("E19").Value = '  +8.03%  '

$ws.Range("D20").Value = '''7.75'
$ws.Range("E20").Value = '  +11.20%  '

$ws.Range("D21").Value = '''324.09'
$ws.Range("E21").Value = '  +6.43%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("E23").Value = '  +6.33%  '

$ws.Range("D24").Value = '''61.92'
$ws.Range("E24").Value = '  +5.19%  '

$ws.Range("D25").Value = '''0.996'
$ws.Range("E25").Value = '  -0.69%  '

$ws.Range("E26").Value = '  +5.40%  '

$ws.Range("D27").Value = '0.0₃0890'
$ws.Range("E27").Value = '  +12.88%  '

$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("E29").Value = '  +11.95%  '

$ws.Range("D30").Value = '''6.77'
$ws.Range("E30").Value = '  +13.51%  '

$ws.Range("E31").Value = '  +7.36%  '

$ws.Range("D32").Value = '''1.75'
$ws.Range("E32").Value = '  +10.95%  '

$ws.Range("D33").Value = '''20.43'
$ws.Range("E33").Value = '  +8.27%  '

$ws.Range("D34").Value = '''157.62'
$ws.Range("E34").Value = '  +9.46%  '

$ws.Range("D35").Value = '''4.45'
$ws.Range("E35").Value = '  +7.30%  '

$ws.Range("E36").Value = '  +5.46%  '

$ws.Range("D37").Value = '''5.54'
$ws.Range("E37").Value = '  +3.89%  '

$ws.Range("D38").Value = '''0.0673'
$ws.Range("E38").Value = '  +11.84%  '

$ws.Range("D39").Value = '''22.97'
$ws.Range("E39").Value = '  +8.42%  '

$ws.Range("D40").Value = '3.015.52'
$ws.Range("E40").Value = '  +2.40%  '

$ws.Range("D41").Value = '''0.998'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("D42").Value = '''36.24'
$ws.Range("E42").Value = '  +3.87%  '

$ws.Range("E43").Value = '  +7.68%  '

$ws.Range("D44").Value = '2.247.72'
$ws.Range("E44").Value = '  +10.96%  '

$ws.Range("D45").Value = '''0.987'
$ws.Range("E45").Value = '  +3.93%  '

$ws.Range("E46").Value = '  +8.64%  '

$ws.Range("D47").Value = '''3.56'
$ws.Range("E47").Value = '  +5.74%  '

$ws.Range("D48").Value = '''1.93'
$ws.Range("E48").Value = '  +26.13%  '

$ws.Range("D49").Value = '''0.0235'
$ws.Range("E49").Value = '  +13.22%  '

$ws.Range("D50").Value = '''5.76'
$ws.Range("E50").Value = '  +10.13%  '

$ws.Range("D51").Value = '''18.92'
$ws.Range("E51").Value = '  +9.10%  '
